$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 551 ("Fruta / hortaliza, semanal").
# This pushes all the existing records that were in rows 551-670 down by one row
# (to rows 552-671), and the freshly inserted row 551 gets the new observation.
$ws.Rows.Item(551).Insert()

$ws.Range("A551").Value = 3
$ws.Range("B551").Value = "Femacal de La Calera"
$ws.Range("C551").Value = "Coquimbo"
$ws.Range("D551").Value = 45244
$ws.Range("E551").Value = 5
$ws.Range("F551").Value = 100112040
$ws.Range("G551").Value = "Cilantro"
$ws.Range("H551").Value = "Sin especificar"
$ws.Range("I551").Value = "Primera"
$ws.Range("J551").Value = 50
$ws.Range("K551").Value = 5000
$ws.Range("L551").Value = 5000
$ws.Range("M551").Value = 5000
$ws.Range("N551").Value = '$/docena de atados (3 kilos)'
$ws.Range("O551").Value = "Provincia de Quillota"
$ws.Range("P551").Value = 1667
$ws.Range("Q551").Value = 3
$ws.Range("R551").Value = "Hortaliza"
